$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the date-formatted style (numFmtId 14) already used by A2/A3
# instead of letting NumberFormat mint a brand-new custom format.
$ws.Range("A2").Copy()

# Row 4 - week of 2017-10-06 (serial 43014)
$ws.Range("A4").PasteSpecial(-4122)
$ws.Cells.Item(4, 1).Value = 43014
$ws.Cells.Item(4, 2).Value = "Android Main GUI first draft; Pupil and Database class"
$ws.Cells.Item(4, 3).Value = "C# Registration: checks if values are correct"
$ws.Cells.Item(4, 4).Value = "Research MongoDB"

# Row 5 - week of 2017-10-13 (serial 43021)
$ws.Range("A5").PasteSpecial(-4122)
$ws.Cells.Item(5, 1).Value = 43021
$ws.Cells.Item(5, 2).Value = "Android Login bugs fixed"
$ws.Cells.Item(5, 3).Value = "C# MainWindow with Menu finished; AddEntryWindow finished"
$ws.Cells.Item(5, 4).Value = "Research MongoDB (not finished), Virtual Machine (Aphrodite? Not finished)"

# Row 6 - week of 2017-10-20 (serial 43028)
$ws.Range("A6").PasteSpecial(-4122)
$ws.Cells.Item(6, 1).Value = 43028
$ws.Cells.Item(6, 3).Value = "C# AddEntry: checks if values are correct"
$ws.Cells.Item(6, 2).Value = "Android AddEntryWindow finished"
$ws.Cells.Item(6, 4).Value = "Virtual Machine bug fixed. MongoDB locally installed; Table company created"

$ws.Range("D6").Select()
